$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new draw result as row 56, matching the existing rows: every
# cell is plain text (dates/phase codes/results are stored as strings, not
# numbers or Excel dates). Force text typing via a temporary "@" number
# format so values like "2025-11-11" / "251111" aren't auto-coerced into a
# date serial / number, then clear the format again so the new cells end up
# on the same (default) style as the rest of the sheet.
$row = 56

$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-11-11"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = "Pick 3"

$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value = "251111"
$ws.Range("C$row").ClearFormats()

$ws.Range("D$row").Value = "2-2-2"

$ws.Range("E$row").Value = "2025-11-11T21:40:24.499+04:00"
